{"js": "// Commit: \"Use a full row retrieval method\"\n//\n// The `{% for blurb in pet[p].blurbs %}` loop is changed to iterate over\n// full dataframe rows instead of a pre-built \"blurbs\" list:\n//   {%p for blurb in pet[p].blurbs %}      ->  {%p for row in pet[p].selected_labels_df.iterrows() %}\n// and the body reference is updated to pull the \"Blurb\" column out of the\n// row tuple returned by `iterrows()`:\n//   {{ blurb }}                             ->  {{ row[1][\"Blurb\"] }}\n\n// 1) Rename the loop variable \"blurb\" -> \"row\" everywhere it is used as a\n//    standalone identifier (the \"for X in ...\" variable and its use in the\n//    body), without touching the \"blurbs\" attribute name.\nconst loopVarMatches = context.document.body.search(\"blurb\", {\n  matchCase: true,\n  matchWholeWord: true,\n});\nloopVarMatches.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < loopVarMatches.items.length; i++) {\n  loopVarMatches.items[i].insertText(\"row\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Swap the \"blurbs\" list attribute access for a full-row retrieval call:\n//    pet[p].blurbs  ->  pet[p].selected_labels_df.iterrows()\nconst attrMatches = context.document.body.search(\"blurbs \", {\n  matchCase: true,\n});\nattrMatches.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < attrMatches.items.length; i++) {\n  attrMatches.items[i].insertText(\n    \"selected_labels_df.iterrows() \",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// 3) Update the display expression to index into the row tuple and pull out\n//    the \"Blurb\" column: {{ row }} -> {{ row[1][\"Blurb\"] }}\nconst displayMatches = context.document.body.search(\"{{ row }}\", {\n  matchCase: true,\n});\ndisplayMatches.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < displayMatches.items.length; i++) {\n  displayMatches.items[i].insertText(\n    '{{ row[1][\"Blurb\"] }}',\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n", "ps1": "# Commit: \"Use a full row retrieval method\"\n#\n# The `{% for blurb in pet[p].blurbs %}` loop is changed to iterate over\n# full dataframe rows instead of a pre-built \"blurbs\" list:\n#   {%p for blurb in pet[p].blurbs %}      ->  {%p for row in pet[p].selected_labels_df.iterrows() %}\n# and the body reference is updated to pull the \"Blurb\" column out of the\n# row tuple returned by `iterrows()`:\n#   {{ blurb }}                             ->  {{ row[1][\"Blurb\"] }}\n#\n# NOTE: replacement text is written via Range.Text (rather than\n# Find.Execute's ReplaceWith parameter) so the straight double-quotes in\n# `[\"Blurb\"]` are not auto-corrected into curly/smart quotes.\n\n$d = $word.ActiveDocument\n\n# 1) Rename the loop variable \"blurb\" -> \"row\" everywhere it is used as a\n#    standalone identifier (the \"for X in ...\" variable and its use in the\n#    body), without touching the \"blurbs\" attribute name.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"blurb\"\n$find.MatchWholeWord = $true\n$find.MatchCase = $true\n$find.Forward = $true\n$find.Wrap = 0  # wdFindStop\nwhile ($find.Execute()) {\n    $rng.Text = \"row\"\n    $rng.Collapse(0)  # wdCollapseEnd\n    $rng.End = $d.Content.End\n}\n\n# 2) Swap the \"blurbs\" list attribute access for a full-row retrieval call:\n#    pet[p].blurbs  ->  pet[p].selected_labels_df.iterrows()\n$rng2 = $d.Content\n$find2 = $rng2.Find\n$find2.ClearFormatting()\n$find2.Text = \"blurbs \"\n$find2.MatchWholeWord = $false\n$find2.MatchCase = $true\n$find2.Forward = $true\n$find2.Wrap = 0  # wdFindStop\nif ($find2.Execute()) {\n    $rng2.Text = \"selected_labels_df.iterrows() \"\n}\n\n# 3) Update the display expression to index into the row tuple and pull out\n#    the \"Blurb\" column: {{ row }} -> {{ row[1][\"Blurb\"] }}\n$rng3 = $d.Content\n$find3 = $rng3.Find\n$find3.ClearFormatting()\n$find3.Text = \"{{ row }}\"\n$find3.MatchWholeWord = $false\n$find3.MatchCase = $true\n$find3.Forward = $true\n$find3.Wrap = 0  # wdFindStop\nif ($find3.Execute()) {\n    $rng3.Text = '{{ row[1][\"Blurb\"] }}'\n}\n"}
